$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "ODI Batting" sheet: a handful of rows have an explicitly-blank
#    INNING_NUMBER cell (B) for matches the player did not bat in.
#    Those placeholder cells get dropped entirely.
# ------------------------------------------------------------------
$battingWs = $wb.Worksheets.Item("ODI Batting")
$battingWs.Range("B2").ClearContents()
$battingWs.Range("B4").ClearContents()
$battingWs.Range("B5").ClearContents()
$battingWs.Range("B6").ClearContents()
$battingWs.Range("B8").ClearContents()
$battingWs.Range("B9").ClearContents()

# ------------------------------------------------------------------
# 2) Add a new trailing sheet "ODI Batting Extra" with extra batting
#    stats per match.
# ------------------------------------------------------------------
$headerSrc = $battingWs.Range("A1:F1")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Header row text first ...
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# ... then stamp it with the same bold/centered/bordered header style
# used by the other sheets (reuse the style, don't fabricate a new one).
$headerSrc.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# Data rows. MATCH_CODE (A) is numeric-looking text, so it is entered
# with a leading apostrophe to keep it as text instead of a number.
# Same for the NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL columns, which
# hold text such as "0" or "1.83%" rather than real numbers. Where the
# source data has no value at all the cell is still written as an
# (empty) text cell via a bare apostrophe, matching the source export.
$ws.Range("A2").Value = "'4234"
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = "'"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'"
$ws.Range("F2").Value = "NO"

$ws.Range("A3").Value = "'4235"
$ws.Range("B3").Value = "'"
$ws.Range("C3").Value = "'"
$ws.Range("D3").Value = "'"
$ws.Range("E3").Value = "'"
$ws.Range("F3").Value = "NO"

$ws.Range("A4").Value = "'4258"
$ws.Range("B4").Value = "'"
$ws.Range("C4").Value = "'"
$ws.Range("D4").Value = "'"
$ws.Range("E4").Value = "'"
$ws.Range("F4").Value = "NO"

$ws.Range("A5").Value = "'4268"
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = "'"
$ws.Range("D5").Value = "'"
$ws.Range("E5").Value = "'"
$ws.Range("F5").Value = "NO"

$ws.Range("A6").Value = "'4275"
$ws.Range("B6").Value = 11
$ws.Range("C6").Value = "'"
$ws.Range("D6").Value = "'"
$ws.Range("E6").Value = "'"
$ws.Range("F6").Value = "NO"

$ws.Range("A7").Value = "'4277"
$ws.Range("B7").Value = 8
$ws.Range("C7").Value = "'0"
$ws.Range("D7").Value = "'0"
$ws.Range("E7").Value = "'1.83%"
$ws.Range("F7").Value = "NO"

$ws.Range("A8").Value = "'4322"
$ws.Range("B8").Value = 10
$ws.Range("C8").Value = "'"
$ws.Range("D8").Value = "'"
$ws.Range("E8").Value = "'"
$ws.Range("F8").Value = "NO"

$ws.Range("A9").Value = "'4336"
$ws.Range("B9").Value = 11
$ws.Range("C9").Value = "'"
$ws.Range("D9").Value = "'"
$ws.Range("E9").Value = "'"
$ws.Range("F9").Value = "NO"

$ws.Range("A10").Value = "'4341"
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = "'0"
$ws.Range("D10").Value = "'0"
$ws.Range("E10").Value = "'"
$ws.Range("F10").Value = "NO"

$ws.Range("A11").Value = "'4351"
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = "'0"
$ws.Range("D11").Value = "'0"
$ws.Range("E11").Value = "'3.49%"
$ws.Range("F11").Value = "NO"

$ws.Range("A12").Value = "'4354"
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = "'0"
$ws.Range("D12").Value = "'0"
$ws.Range("E12").Value = "'0.45%"
$ws.Range("F12").Value = "NO"

$ws.Range("A13").Value = "'4567"
$ws.Range("B13").Value = 9
$ws.Range("C13").Value = "'0"
$ws.Range("D13").Value = "'0"
$ws.Range("E13").Value = "'0.95%"
$ws.Range("F13").Value = "NO"
